$d = $word.ActiveDocument

# Locate the paragraph that contains "Teste 1 ... alterando documento.."
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Teste 1*alterando documento*") {
        $target = $p
    }
}

if ($target -eq $null) {
    Write-Host "Target paragraph not found"
} else {
    # Build a zero-length range positioned one character before the end of
    # the paragraph's text (i.e. just after the last visible character but
    # before the paragraph mark). Using the exact paragraph-mark boundary
    # position causes InsertXML to swallow/merge an adjoining paragraph
    # mark, so we back up by one character to stay safely inside the run.
    $insertAt = $target.Range.End - 1
    $r = $d.Range($insertAt, $insertAt)

    $xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<w:wordDocument xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>
</w:p>
<w:p>
<w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>
<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Teste 2</w:t></w:r>
</w:p>
</w:body>
</w:wordDocument>
'@

    [void]$r.InsertXML($xml)
}
